# Update the "Xml" (column C) benchmark values on the "Courses" sheet
# with the re-measured timings, then restore the "Courses" sheet as the
# active / selected sheet (it had lost focus to "Prices").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Courses")

$values = @(
    5.9041160000000001,
    6.4622999999999999,
    7.7937500000000002,
    8.7020999999999997,
    9.1550829999999994,
    12.400766000000001,
    22.224381999999999,
    37.048099999999998,
    54.853316,
    139.180566
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# Make "Courses" the active sheet again and select C2:C11 on it, with
# the active cell at the top of the range (C2).
$ws.Activate()
$ws.Range("C2:C11").Select()
